{"js": "// Word Diario di bordo: merge the \"29/1/25\" date run fragments, then\n// append a new diary entry (date heading \"5/2/25\" + body paragraph)\n// right after the \"Aggiunta codice...\" entry, describing the completed\n// AMC code and successful tests.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Locate the two anchor paragraphs by their text content rather than by\n// a fixed index, so the script is resilient to unrelated document\n// changes.\nparagraphs.items.forEach((p) => p.load(\"text,style\"));\nawait context.sync();\n\nlet dateParagraph = null;\nlet kmeansParagraph = null;\nfor (const p of paragraphs.items) {\n  const text = p.text.trim();\n  if (text === \"29/1/25\" || text === \"29\" + \"/1/2\" + \"5\") {\n    dateParagraph = p;\n  }\n  if (text.indexOf(\"Aggiunta codice per creazione del dataset di training\") === 0) {\n    kmeansParagraph = p;\n  }\n}\n\nif (!dateParagraph || !kmeansParagraph) {\n  throw new Error(\"Could not locate anchor paragraphs for the 29/1/25 entry.\");\n}\n\n// Re-write the date heading so its text ends up as a single run\n// (\"29\" + \"/1/2\" + \"5\" -> \"29/1/25\").\ndateParagraph.clear();\ndateParagraph.insertText(\"29/1/25\", Word.InsertLocation.start);\n\n// Insert the new \"5/2/25\" heading paragraph right after the K-means\n// entry, matching the Heading 2 style used by the other date markers.\nconst newDateParagraph = kmeansParagraph.insertParagraph(\"5/2/25\", Word.InsertLocation.after);\nnewDateParagraph.styleBuiltIn = Word.Style.heading2;\n\n// Insert the new diary body paragraph right after the new date heading.\nconst newBodyParagraph = newDateParagraph.insertParagraph(\n  \"Il codice per l\\u2019AMC \\u00e8 concluso. I test sono stati un successo. Bisogna solo registrare delle prove con un canale con ostacoli in movimento. Si pensa ora ad andare oltre. Bisogna rifare le misurazioni per le codifiche.\",\n  Word.InsertLocation.after\n);\nnewBodyParagraph.styleBuiltIn = Word.Style.normal;\n\nawait context.sync();\n", "ps1": "# Word Diario di bordo: merge the \"29/1/25\" date run fragments into a\n# single run, then append a new diary entry (date heading \"5/2/25\" +\n# body paragraph) right after the \"Aggiunta codice...\" entry, describing\n# the completed AMC code and successful tests.\n\n$d = $word.ActiveDocument\n\n# --- Locate the two anchor paragraphs by their text content rather than\n# a fixed index, so the script is resilient to unrelated document\n# changes. ---\n$dateParagraph = $null\n$kmeansParagraph = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text.Trim()\n    if ($t -eq \"29/1/25\") {\n        $dateParagraph = $p\n    }\n    if ($t.StartsWith(\"Aggiunta codice per creazione del dataset di training\")) {\n        $kmeansParagraph = $p\n    }\n}\n\nif ($dateParagraph -eq $null -or $kmeansParagraph -eq $null) {\n    throw \"Could not locate anchor paragraphs for the 29/1/25 entry.\"\n}\n\n# --- Re-write the date heading so its text ends up as a single run\n# (\"29\" + \"/1/2\" + \"5\" -> \"29/1/25\"). Using Find/Replace on the\n# paragraph's own range cleanly substitutes the whole (multi-run) match\n# with one new run. ---\n$dateRange = $dateParagraph.Range\n$find = $dateRange.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"29/1/25\", $false, $false, $false, $false, $false, $true, 1, $false, \"29/1/25\", 2) | Out-Null\n\n# --- Insert a new paragraph right after the K-means entry, set its text\n# to the new date heading, and apply the Heading 2 style used by the\n# other date markers. ---\n$kmeansParagraph.Range.InsertParagraphAfter()\n$newDateParagraph = $d.Paragraphs.Item($kmeansParagraph.Index + 1)\n$newDateParagraph.Range.Text = \"5/2/25\"\n$newDateParagraph.Style = \"Heading 2\"\n\n# --- Insert the new diary body paragraph right after the new date\n# heading. ---\n$newDateParagraph.Range.InsertParagraphAfter()\n$newBodyParagraph = $d.Paragraphs.Item($newDateParagraph.Index + 1)\n$newBodyParagraph.Range.Text = \"Il codice per l\u2019AMC \u00e8 concluso. I test sono stati un successo. Bisogna solo registrare delle prove con un canale con ostacoli in movimento. Si pensa ora ad andare oltre. Bisogna rifare le misurazioni per le codifiche.\"\n$newBodyParagraph.Style = \"Normal\"\n"}
